$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1649484536082474
$ws.Range("C2").Value = 0.6254295532646048
$ws.Range("J2").Value = 0.0274914089347079
$ws.Range("P2").Value = 0.1237113402061856
$ws.Range("S2").Value = 0.0584192439862543
$ws.Range("B3").Value = 0.005263157894736842
$ws.Range("C3").Value = 0.01578947368421053
$ws.Range("J3").Value = 0.04736842105263158
$ws.Range("P3").Value = 0.8052631578947368
$ws.Range("S3").Value = 0.1263157894736842
$ws.Range("J4").Value = 0.04347826086956522
$ws.Range("P4").Value = 0.7536231884057971
$ws.Range("S4").Value = 0.2028985507246377
$ws.Range("B6").Value = 0.04405286343612335
$ws.Range("D6").Value = 0.02643171806167401
$ws.Range("F6").Value = 0.05286343612334802
$ws.Range("J6").Value = 0.2775330396475771
$ws.Range("O6").Value = 0.01762114537444934
$ws.Range("Q6").Value = 0.1894273127753304
$ws.Range("R6").Value = 0.09251101321585903
$ws.Range("S6").Value = 0.2995594713656388
$ws.Range("B7").Value = 0.07106598984771574
$ws.Range("D7").Value = 0.06091370558375635
$ws.Range("F7").Value = 0.07106598984771574
$ws.Range("J7").Value = 0.1472081218274112
$ws.Range("O7").Value = 0.02538071065989848
$ws.Range("Q7").Value = 0.2284263959390863
$ws.Range("R7").Value = 0.08121827411167512
$ws.Range("S7").Value = 0.3147208121827411
$ws.Range("B8").Value = 0.08955223880597014
$ws.Range("D8").Value = 0.02345415778251599
$ws.Range("F8").Value = 0.07462686567164178
$ws.Range("J8").Value = 0.1343283582089552
$ws.Range("O8").Value = 0.02345415778251599
$ws.Range("Q8").Value = 0.1684434968017058
$ws.Range("R8").Value = 0.1044776119402985
$ws.Range("S8").Value = 0.3816631130063966
$ws.Range("B9").Value = 0.08484848484848485
$ws.Range("D9").Value = 0.01818181818181818
$ws.Range("F9").Value = 0.08484848484848485
$ws.Range("J9").Value = 0.1393939393939394
$ws.Range("O9").Value = 0.01818181818181818
$ws.Range("Q9").Value = 0.2666666666666667
$ws.Range("R9").Value = 0.06060606060606061
$ws.Range("S9").Value = 0.3272727272727273
$ws.Range("B10").Value = 0.1053704962610469
$ws.Range("D10").Value = 0.02515295717199184
$ws.Range("F10").Value = 0.06254248810333107
$ws.Range("J10").Value = 0.1250849762066621
$ws.Range("O10").Value = 0.01155676410605031
$ws.Range("Q10").Value = 0.2195785180149558
$ws.Range("R10").Value = 0.09925220938137322
$ws.Range("S10").Value = 0.3514615907545887
$ws.Range("G11").Value = 0.1453287197231834
$ws.Range("J11").Value = 0.08996539792387544
$ws.Range("K11").Value = 0.2110726643598616
$ws.Range("L11").Value = 0.5432525951557093
$ws.Range("S11").Value = 0.01038062283737024
$ws.Range("G12").Value = 0.7423312883435583
$ws.Range("J12").Value = 0.1717791411042945
$ws.Range("K12").Value = 0.006134969325153374
$ws.Range("L12").Value = 0.049079754601227
$ws.Range("S12").Value = 0.03067484662576687
$ws.Range("G13").Value = 0.7551020408163265
$ws.Range("J13").Value = 0.2244897959183673
$ws.Range("S13").Value = 0.02040816326530612
$ws.Range("F15").Value = 0.01968503937007874
$ws.Range("H15").Value = 0.1574803149606299
$ws.Range("I15").Value = 0.05905511811023622
$ws.Range("J15").Value = 0.4251968503937008
$ws.Range("K15").Value = 0.05118110236220472
$ws.Range("M15").Value = 0.007874015748031496
$ws.Range("O15").Value = 0.06299212598425197
$ws.Range("S15").Value = 0.2165354330708661
$ws.Range("F16").Value = 0.03125
$ws.Range("H16").Value = 0.1919642857142857
$ws.Range("I16").Value = 0.0625
$ws.Range("J16").Value = 0.4464285714285715
$ws.Range("K16").Value = 0.09375
$ws.Range("M16").Value = 0.02678571428571428
$ws.Range("N16").Value = 0.004464285714285714
$ws.Range("O16").Value = 0.05803571428571429
$ws.Range("S16").Value = 0.08482142857142858
$ws.Range("F17").Value = 0.0149812734082397
$ws.Range("H17").Value = 0.1760299625468165
$ws.Range("I17").Value = 0.08052434456928839
$ws.Range("J17").Value = 0.4456928838951311
$ws.Range("K17").Value = 0.09925093632958802
$ws.Range("M17").Value = 0.01685393258426966
$ws.Range("O17").Value = 0.05805243445692884
$ws.Range("S17").Value = 0.1086142322097378
$ws.Range("F18").Value = 0.01680672268907563
$ws.Range("H18").Value = 0.2100840336134454
$ws.Range("I18").Value = 0.06722689075630252
$ws.Range("J18").Value = 0.4789915966386555
$ws.Range("K18").Value = 0.07983193277310924
$ws.Range("M18").Value = 0.008403361344537815
$ws.Range("O18").Value = 0.0546218487394958
$ws.Range("S18").Value = 0.08403361344537816
$ws.Range("F19").Value = 0.01313628899835796
$ws.Range("H19").Value = 0.2011494252873563
$ws.Range("I19").Value = 0.06157635467980296
$ws.Range("J19").Value = 0.396551724137931
$ws.Range("K19").Value = 0.09441707717569786
$ws.Range("M19").Value = 0.0270935960591133
$ws.Range("N19").Value = 0.001642036124794745
$ws.Range("O19").Value = 0.08702791461412152
$ws.Range("S19").Value = 0.1174055829228243
